# Notation.xlsx update
# - verified that the rearranged equations could be replicated in code,
#   updated variable names in code to reflect Notation.xlsx.
# - Updated some sections of workbook (Code sheet gets two new
#   python-variable-name entries; Math sheet becomes the active/selected
#   tab again).

$wb = $excel.ActiveWorkbook
$wsMath = $wb.Worksheets.Item("Math")
$wsCode = $wb.Worksheets.Item("Code")

# --- Code sheet content updates ------------------------------------------
# New row: "orthogonal_sum" python variable lines up with the cCon /
# "central contribution offset" row (row 8), in the Python variable column G.
$wsCode.Range("G8").Value = "orthogonal_sum"
$wsCode.Range("G8").WrapText = $true

# Row 10 (second "i" / index row) python variable renamed from
# component_index -> component_array_index.
$wsCode.Range("G10").Value = "component_array_index"

# --- Active tab / selection swap -----------------------------------------
# Previously "Code" was the active/selected tab with selection D5.
# Now "Math" is the active/selected tab with selection B28, and "Code"'s
# stored selection moves to G9.
$wsCode.Range("G9").Select()
$wsMath.Activate()
$wsMath.Range("B28").Select()
